$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two auxiliary AVERAGE() summary rows (119, 120) and the extra
# AVERAGE() row (121) that lived below the raw data table are no longer
# needed now that the vehicles have been re-classified -- remove them
# entirely (bottom-up so row numbers of earlier rows are not disturbed).
[void]$ws.Rows.Item(121).Delete()
[void]$ws.Rows.Item(120).Delete()
[void]$ws.Rows.Item(119).Delete()

# Re-classification is complete: every data row (2-118) should now be
# visible again instead of being hidden behind the old AutoFilter criteria.
for ($r = 2; $r -le 118; $r++) {
  $ws.Rows.Item($r).Hidden = $false
}

# Drop the AutoFilter entirely (criteria + dropdown buttons + the
# _FilterDatabase plumbing it implies).
$ws.AutoFilterMode = $false

# Remove any leftover defined names (e.g. the hidden _xlnm._FilterDatabase
# name Excel maintains for the AutoFilter range).
foreach ($n in @($wb.Names)) {
  [void]$n.Delete()
}

# The reference picture/legend that illustrated the old fit is no longer
# relevant to the new survival-pattern fit -- remove it.
if ($ws.Shapes.Count -gt 0) {
  [void]$ws.Shapes.Item(1).Delete()
}

# Reset the view back to the top-left cell instead of the stray selection
# left over from the filtering/plotting session.
[void]$ws.Range("A1").Select()
